$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43-67 down to 44-68.
$ws.Rows.Item(43).Insert()

# Fill the newly inserted row 43 with the new weekly record.
$ws.Range("A43").Value = 6
$ws.Range("B43").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C43").Value = "Metropolitana"
$ws.Range("D43").Value = 45062
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 100112035
$ws.Range("G43").Value = "Bruselas (repollito)"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 300
$ws.Range("K43").Value = 18000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = 18800
$ws.Range("N43").Value = "$/malla 15 kilos"
$ws.Range("O43").Value = "Provincia de Quillota"
$ws.Range("P43").Value = 1253
$ws.Range("Q43").Value = 15
$ws.Range("R43").Value = "Hortaliza"
